# Reorder the data rows of Sheet1 (rows 2-11, row 1 is the header) so that
# the rows end up in the sequence: 2,7,8,9,10,11,6,3,4,5 (original row numbers).
# Row 1 (header) and row 2 stay put; only rows 3-11 actually move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture the current contents of the 6 data columns (A:F) for rows 2-11
# before we start overwriting anything.
$original = @{}
for ($r = 2; $r -le 11; $r++) {
    $original[$r] = $ws.Range("A" + $r + ":F" + $r).Value2
}

# Mapping of new row number -> original row number it should now contain.
$newOrder = @{
    2  = 2
    3  = 7
    4  = 8
    5  = 9
    6  = 10
    7  = 11
    8  = 6
    9  = 3
    10 = 4
    11 = 5
}

for ($newRow = 2; $newRow -le 11; $newRow++) {
    $srcRow = $newOrder[$newRow]
    $ws.Range("A" + $newRow + ":F" + $newRow).Value2 = $original[$srcRow]
}

# Update the active selection on Sheet1 to match the edited workbook.
$ws.Range("B15").Select()
